$d = $word.ActiveDocument

# Locate the DEWALT 20-volt power tools table: the table whose last
# existing row's first cell reads "DCM849P2" (20V MAX* XR 7 in. Cordless
# Variable-Speed Rotary Polisher Kit), then append a new row for the
# DCCS623 20V Max Cordless Pruning Saw.
$target = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $rows = $t.Rows.Count
    $lastCellText = $t.Cell($rows, 1).Range.Text
    $lastCellText = $lastCellText.Replace([char]13, "").Replace([char]7, "")
    if ($lastCellText -eq "DCM849P2") {
        $target = $t
        break
    }
}

$newRow = $target.Rows.Add()
$idx = $newRow.Index
$target.Cell($idx, 1).Range.Text = "DCCS623"
$target.Cell($idx, 2).Range.Text = "20V Max Cordless Pruning Saw"
